$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.679.74'
$ws.Range('E2').Value = '  +3.83%  '
$ws.Range('D3').Value = '1.911.53'
$ws.Range('E3').Value = '  +1.61%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.43'
$ws.Range('E5').Value = '  +1.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.701'
$ws.Range('E6').Value = '  +2.90%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '45.77'
$ws.Range('E8').Value = '  +5.45%  '
$ws.Range('E9').Value = '  +4.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.51'
$ws.Range('E10').Value = '  +9.73%  '
$ws.Range('E11').Value = '  +2.24%  '
$ws.Range('E12').Value = '  +2.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.52'
$ws.Range('E13').Value = '  +7.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.814'
$ws.Range('E14').Value = '  +6.68%  '
$ws.Range('D15').Value = '2.191.73'
$ws.Range('E15').Value = '  +1.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.12'
$ws.Range('E16').Value = '  +3.60%  '
$ws.Range('D17').Value = '1.914.57'
$ws.Range('E17').Value = '  +2.80%  '
$ws.Range('D18').Value = '36.676.56'
$ws.Range('E18').Value = '  +3.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.74'
$ws.Range('E19').Value = '  +1.50%  '
$ws.Range('E20').Value = '  +4.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '250.01'
$ws.Range('E21').Value = '  +2.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.36'
$ws.Range('E22').Value = '  +4.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.19'
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('E24').Value = '  -1.80%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.18'
$ws.Range('E26').Value = '  +0.88%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.80'
$ws.Range('E27').Value = '  +1.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.77'
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.69'
$ws.Range('E29').Value = '  +2.11%  '
$ws.Range('E30').Value = '  +1.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.55'
$ws.Range('E31').Value = '  +5.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0616'
$ws.Range('E32').Value = '  +3.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.32'
$ws.Range('E33').Value = '  +3.17%  '
$ws.Range('E34').Value = '  +22.86%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.87'
$ws.Range('E36').Value = '  +3.07%  '
$ws.Range('E37').Value = '  +5.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.875'
$ws.Range('E38').Value = '  +2.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.74'
$ws.Range('E39').Value = '  +48.30%  '
$ws.Range('E40').Value = '  +2.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '105.16'
$ws.Range('E41').Value = '  +8.26%  '
$ws.Range('E42').Value = '  +3.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.42'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.91'
$ws.Range('E44').Value = '  +21.74%  '
$ws.Range('E45').Value = '  +2.46%  '
$ws.Range('D46').Value = '1.347.71'
$ws.Range('E46').Value = '  +2.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.37'
$ws.Range('E47').Value = '  -1.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0814'
$ws.Range('E48').Value = '  +1.62%  '
$ws.Range('E49').Value = '  +2.46%  '
$ws.Range('E50').Value = '  +2.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.36'
$ws.Range('E51').Value = '  +2.84%  '
